# Add season record columns (Wins, Losses, Ties) to the player stats sheet.
# This mirrors the author's fix: previously only team statistics were pulled,
# not the season win/loss/tie record, so three new columns are appended
# after the existing "Unnamed: 28" column (AC) -> AD:Wins, AE:Losses, AF:Ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold font, border, centered alignment) from an
# existing header cell so the new header cells share the same style.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (rows 2-51): every team had the same 92-70-0 record for 2013 ---
$lastRow = 51
$ws.Range("AD2:AD" + $lastRow).Value = 92
$ws.Range("AE2:AE" + $lastRow).Value = 70
$ws.Range("AF2:AF" + $lastRow).Value = 0

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
